# Re-process the metadata with the newly curated dimensions: sector-descripcion,
# mes-nombre and sexo move from "dimension" to "measure", so their metadata rows
# (type, datatype, mapping file) are updated/cleared accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (iaest-* qualifier): dimension -> measure for the three curated columns
$ws.Range("E2").Value = "iaest-measure:sector-descripcion"
$ws.Range("I2").Value = "iaest-measure:mes-nombre"
$ws.Range("L2").Value = "iaest-measure:sexo"

# Row 3 (dim/medida classifier): dim -> medida
$ws.Range("E3").Value = "medida"
$ws.Range("I3").Value = "medida"
$ws.Range("L3").Value = "medida"

# Row 4 (datatype): skos:Concept -> xsd:int
$ws.Range("E4").Value = "xsd:int"
$ws.Range("I4").Value = "xsd:int"
$ws.Range("L4").Value = "xsd:int"

# Row 5 (mapping file): no longer needed now that these are measures, not
# curated dimensions with a lookup table -> clear the cells.
$ws.Range("E5").Clear()
$ws.Range("I5").Clear()
$ws.Range("L5").Clear()
